# Generate Report for Archive
#
# The localization status report is regenerated: the record for
# "e596abf4-f986-4eba-87ad-410bf4d98f69" now sorts ahead of the record for
# "a5a6e20e-0fcb-4834-b73d-ac720a9a06ff" (it moved from row 7 to row 6, and
# the a5a6e20e record moved from row 6 to row 7) on every sheet, and the
# e596abf4 record's status flips from "Ready for handoff" to
# "In Translation" now that it is back in translation.
#
# Only the cells whose content actually differs between the two records are
# touched, so columns that already hold identical values for both records
# (e.g. booleans, blank placeholders) are left completely untouched.

$wb = $excel.ActiveWorkbook

function Swap-HyperlinkDisplay($ws, $addr6, $addr7) {
    $h6 = $null
    $h7 = $null
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq $addr6) { $h6 = $h }
        if ($addr -eq $addr7) { $h7 = $h }
    }
    if (($h6 -ne $null) -and ($h7 -ne $null)) {
        $t6 = $h6.TextToDisplay
        $t7 = $h7.TextToDisplay
        $h6.TextToDisplay = $t7
        $h7.TextToDisplay = $t6
    }
}

# ---- Overview sheet ----
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A6").Value = "e596abf4-f986-4eba-87ad-410bf4d98f69.md"
$ws1.Range("B6").Value = "e2e\e596abf4-f986-4eba-87ad-410bf4d98f69.md"
$ws1.Range("E6").Value = "In Translation"
$ws1.Range("F6").Value = "In Translation"
$ws1.Range("G6").Value = "2017-02-09 14:50:28"
$ws1.Range("A7").Value = "a5a6e20e-0fcb-4834-b73d-ac720a9a06ff.md"
$ws1.Range("B7").Value = "e2e\a5a6e20e-0fcb-4834-b73d-ac720a9a06ff.md"
$ws1.Range("G7").Value = "2017-02-09 14:44:42"
Swap-HyperlinkDisplay $ws1 '$B$6' '$B$7'

# ---- zh-cn sheet ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A6").Value = "e596abf4-f986-4eba-87ad-410bf4d98f69.md"
$ws2.Range("C6").Value = "In Translation"
$ws2.Range("G6").Value = "e596abf4-f986-4eba-87ad-410bf4d98f69.50bc23ae9beb6d1cc1ac380d58c6a61c9e467441.zh-cn.xlf"
$ws2.Range("H6").Value = "2017-02-09 14:50:09"
$ws2.Range("A7").Value = "a5a6e20e-0fcb-4834-b73d-ac720a9a06ff.md"
$ws2.Range("G7").Value = "a5a6e20e-0fcb-4834-b73d-ac720a9a06ff.a702e7c88a4951d07ec926c96f486c82293c4619.zh-cn.xlf"
$ws2.Range("H7").Value = "2017-02-09 14:44:25"
Swap-HyperlinkDisplay $ws2 '$A$6' '$A$7'

# ---- de-de sheet ----
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A6").Value = "e596abf4-f986-4eba-87ad-410bf4d98f69.md"
$ws3.Range("C6").Value = "In Translation"
$ws3.Range("G6").Value = "e596abf4-f986-4eba-87ad-410bf4d98f69.50bc23ae9beb6d1cc1ac380d58c6a61c9e467441.de-de.xlf"
$ws3.Range("H6").Value = "2017-02-09 14:50:28"
$ws3.Range("A7").Value = "a5a6e20e-0fcb-4834-b73d-ac720a9a06ff.md"
$ws3.Range("G7").Value = "a5a6e20e-0fcb-4834-b73d-ac720a9a06ff.a702e7c88a4951d07ec926c96f486c82293c4619.de-de.xlf"
$ws3.Range("H7").Value = "2017-02-09 14:44:42"
Swap-HyperlinkDisplay $ws3 '$A$6' '$A$7'
